$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 3).Value = 2
$ws.Cells.Item(2, 6).Value = 68
$ws.Cells.Item(2, 12).Value = "stimuli/img_uy1n4.png"
$ws.Cells.Item(2, 13).Value = 76.30555555555556
$ws.Cells.Item(2, 14).Value = 55.33333333333334
$ws.Cells.Item(2, 15).Value = 65.81944444444444
$ws.Cells.Item(2, 16).Value = 36
$ws.Cells.Item(2, 17).Value = 7
$ws.Cells.Item(2, 18).Value = 7
$ws.Cells.Item(2, 19).Value = 7

# Row 3
$ws.Cells.Item(3, 3).Value = 2
$ws.Cells.Item(3, 6).Value = 69
$ws.Cells.Item(3, 12).Value = "stimuli/img_xti0z.png"
$ws.Cells.Item(3, 13).Value = 81.40625
$ws.Cells.Item(3, 14).Value = 61.4375
$ws.Cells.Item(3, 15).Value = 71.421875
$ws.Cells.Item(3, 16).Value = 32
$ws.Cells.Item(3, 17).Value = 8
$ws.Cells.Item(3, 18).Value = 8
$ws.Cells.Item(3, 19).Value = 8

# Row 4
$ws.Cells.Item(4, 3).Value = 2
$ws.Cells.Item(4, 6).Value = 70
$ws.Cells.Item(4, 12).Value = "stimuli/img_p659z.png"
$ws.Cells.Item(4, 13).Value = 84.21621621621621
$ws.Cells.Item(4, 14).Value = 65.37837837837837
$ws.Cells.Item(4, 15).Value = 74.79729729729729
$ws.Cells.Item(4, 16).Value = 37
$ws.Cells.Item(4, 17).Value = 9
$ws.Cells.Item(4, 18).Value = 9
$ws.Cells.Item(4, 19).Value = 9

# Row 5
$ws.Cells.Item(5, 3).Value = 2
$ws.Cells.Item(5, 6).Value = 71
$ws.Cells.Item(5, 12).Value = "stimuli/img_bwo9g.png"
$ws.Cells.Item(5, 13).Value = 64.81818181818181
$ws.Cells.Item(5, 14).Value = 42.36363636363637
$ws.Cells.Item(5, 15).Value = 53.59090909090909
$ws.Cells.Item(5, 17).Value = 4
$ws.Cells.Item(5, 18).Value = 4
$ws.Cells.Item(5, 19).Value = 4

# Row 6
$ws.Cells.Item(6, 3).Value = 2
$ws.Cells.Item(6, 6).Value = 72
$ws.Cells.Item(6, 12).Value = "stimuli/img_cv9qj.png"
$ws.Cells.Item(6, 13).Value = 60.34375
$ws.Cells.Item(6, 14).Value = 35.34375
$ws.Cells.Item(6, 15).Value = 47.84375
$ws.Cells.Item(6, 16).Value = 32
$ws.Cells.Item(6, 17).Value = 3
$ws.Cells.Item(6, 18).Value = 3
$ws.Cells.Item(6, 19).Value = 3

# Row 7
$ws.Cells.Item(7, 3).Value = 2
$ws.Cells.Item(7, 6).Value = 73
$ws.Cells.Item(7, 12).Value = "stimuli/img_i2k07.png"
$ws.Cells.Item(7, 13).Value = 64.25925925925925
$ws.Cells.Item(7, 14).Value = 40.92592592592592
$ws.Cells.Item(7, 15).Value = 52.59259259259259
$ws.Cells.Item(7, 16).Value = 27
$ws.Cells.Item(7, 17).Value = 3
$ws.Cells.Item(7, 18).Value = 3
$ws.Cells.Item(7, 19).Value = 3

# Row 8
$ws.Cells.Item(8, 3).Value = 2
$ws.Cells.Item(8, 6).Value = 74
$ws.Cells.Item(8, 12).Value = "stimuli/img_ncr40.png"
$ws.Cells.Item(8, 13).Value = 75.66666666666667
$ws.Cells.Item(8, 14).Value = 54.27272727272727
$ws.Cells.Item(8, 15).Value = 64.96969696969697
$ws.Cells.Item(8, 16).Value = 33
$ws.Cells.Item(8, 17).Value = 6
$ws.Cells.Item(8, 18).Value = 6
$ws.Cells.Item(8, 19).Value = 6

# Row 9
$ws.Cells.Item(9, 3).Value = 2
$ws.Cells.Item(9, 6).Value = 75
$ws.Cells.Item(9, 12).Value = "stimuli/img_c0me7.png"
$ws.Cells.Item(9, 13).Value = 68.4
$ws.Cells.Item(9, 14).Value = 45.62857142857143
$ws.Cells.Item(9, 15).Value = 57.01428571428572
$ws.Cells.Item(9, 16).Value = 35
$ws.Cells.Item(9, 17).Value = 4
$ws.Cells.Item(9, 18).Value = 4
$ws.Cells.Item(9, 19).Value = 4

# Row 10
$ws.Cells.Item(10, 3).Value = 2
$ws.Cells.Item(10, 6).Value = 76
$ws.Cells.Item(10, 8).Value = "kitchens"
$ws.Cells.Item(10, 9).Value = "target"
$ws.Cells.Item(10, 11).Value = "j"
$ws.Cells.Item(10, 12).Value = "stimuli/img_ifebc.png"
$ws.Cells.Item(10, 13).Value = 84.0
$ws.Cells.Item(10, 14).Value = 65.88235294117646
$ws.Cells.Item(10, 15).Value = 74.94117647058823
$ws.Cells.Item(10, 16).Value = 34
$ws.Cells.Item(10, 18).Value = 9
$ws.Cells.Item(10, 19).Value = 9

# Row 11
$ws.Cells.Item(11, 3).Value = 2
$ws.Cells.Item(11, 6).Value = 77
$ws.Cells.Item(11, 12).Value = "stimuli/img_ua9bs.png"
$ws.Cells.Item(11, 13).Value = 82.0
$ws.Cells.Item(11, 14).Value = 62.23333333333333
$ws.Cells.Item(11, 15).Value = 72.11666666666667
$ws.Cells.Item(11, 16).Value = 30
$ws.Cells.Item(11, 17).Value = 9
$ws.Cells.Item(11, 18).Value = 9
$ws.Cells.Item(11, 19).Value = 9

# Row 12
$ws.Cells.Item(12, 3).Value = 2
$ws.Cells.Item(12, 6).Value = 78
$ws.Cells.Item(12, 8).Value = "bedrooms"
$ws.Cells.Item(12, 12).Value = "stimuli/img_bklr1.png"
$ws.Cells.Item(12, 13).Value = 86.54761904761905
$ws.Cells.Item(12, 14).Value = 67.73809523809524
$ws.Cells.Item(12, 15).Value = 77.14285714285714
$ws.Cells.Item(12, 16).Value = 42
$ws.Cells.Item(12, 17).Value = 9
$ws.Cells.Item(12, 18).Value = 9
$ws.Cells.Item(12, 19).Value = 9

# Row 13
$ws.Cells.Item(13, 3).Value = 2
$ws.Cells.Item(13, 6).Value = 79
$ws.Cells.Item(13, 8).Value = "kitchens"
$ws.Cells.Item(13, 9).Value = "target"
$ws.Cells.Item(13, 11).Value = "j"
$ws.Cells.Item(13, 12).Value = "stimuli/img_xesl0.png"
$ws.Cells.Item(13, 13).Value = 69.28571428571429
$ws.Cells.Item(13, 14).Value = 47.35714285714285
$ws.Cells.Item(13, 15).Value = 58.32142857142857
$ws.Cells.Item(13, 16).Value = 28
$ws.Cells.Item(13, 17).Value = 5
$ws.Cells.Item(13, 18).Value = 5
$ws.Cells.Item(13, 19).Value = 5

# Row 14
$ws.Cells.Item(14, 3).Value = 2
$ws.Cells.Item(14, 6).Value = 80
$ws.Cells.Item(14, 12).Value = "stimuli/img_j5rpx.png"
$ws.Cells.Item(14, 13).Value = 72.24242424242425
$ws.Cells.Item(14, 14).Value = 50.0
$ws.Cells.Item(14, 15).Value = 61.12121212121212
$ws.Cells.Item(14, 17).Value = 5
$ws.Cells.Item(14, 18).Value = 5
$ws.Cells.Item(14, 19).Value = 5

# Row 15
$ws.Cells.Item(15, 3).Value = 2
$ws.Cells.Item(15, 6).Value = 81
$ws.Cells.Item(15, 8).Value = "living_rooms"
$ws.Cells.Item(15, 9).Value = "distractor"
$ws.Cells.Item(15, 11).Value = "f"
$ws.Cells.Item(15, 12).Value = "stimuli/img_xzyzy.png"
$ws.Cells.Item(15, 13).Value = 85.37209302325581
$ws.Cells.Item(15, 14).Value = 68.90697674418605
$ws.Cells.Item(15, 15).Value = 77.13953488372093
$ws.Cells.Item(15, 16).Value = 43

# Row 16
$ws.Cells.Item(16, 3).Value = 2
$ws.Cells.Item(16, 6).Value = 82
$ws.Cells.Item(16, 12).Value = "stimuli/img_jz3kd.png"
$ws.Cells.Item(16, 13).Value = 72.79411764705883
$ws.Cells.Item(16, 14).Value = 51.64705882352941
$ws.Cells.Item(16, 15).Value = 62.22058823529412
$ws.Cells.Item(16, 16).Value = 34
$ws.Cells.Item(16, 17).Value = 6
$ws.Cells.Item(16, 18).Value = 6
$ws.Cells.Item(16, 19).Value = 6

# Row 17
$ws.Cells.Item(17, 3).Value = 2
$ws.Cells.Item(17, 6).Value = 83
$ws.Cells.Item(17, 12).Value = "stimuli/img_411xa.png"
$ws.Cells.Item(17, 13).Value = 51.03030303030303
$ws.Cells.Item(17, 14).Value = 28.93939393939394
$ws.Cells.Item(17, 15).Value = 39.98484848484848
$ws.Cells.Item(17, 17).Value = 2
$ws.Cells.Item(17, 18).Value = 2
$ws.Cells.Item(17, 19).Value = 2

# Row 18
$ws.Cells.Item(18, 3).Value = 2
$ws.Cells.Item(18, 6).Value = 84
$ws.Cells.Item(18, 12).Value = "stimuli/img_c79r7.png"
$ws.Cells.Item(18, 13).Value = 56.26470588235294
$ws.Cells.Item(18, 14).Value = 34.26470588235294
$ws.Cells.Item(18, 15).Value = 45.26470588235294
$ws.Cells.Item(18, 16).Value = 34
$ws.Cells.Item(18, 17).Value = 2
$ws.Cells.Item(18, 18).Value = 2
$ws.Cells.Item(18, 19).Value = 2

# Row 19
$ws.Cells.Item(19, 3).Value = 2
$ws.Cells.Item(19, 6).Value = 85
$ws.Cells.Item(19, 8).Value = "bedrooms"
$ws.Cells.Item(19, 9).Value = "distractor"
$ws.Cells.Item(19, 11).Value = "f"
$ws.Cells.Item(19, 12).Value = "stimuli/img_ca8kd.png"
$ws.Cells.Item(19, 13).Value = 92.05405405405405
$ws.Cells.Item(19, 14).Value = 73.02702702702703
$ws.Cells.Item(19, 15).Value = 82.54054054054055
$ws.Cells.Item(19, 16).Value = 37
$ws.Cells.Item(19, 17).Value = 10
$ws.Cells.Item(19, 18).Value = 10
$ws.Cells.Item(19, 19).Value = 10

# Row 20
$ws.Cells.Item(20, 3).Value = 2
$ws.Cells.Item(20, 6).Value = 86
$ws.Cells.Item(20, 12).Value = "stimuli/img_i7vab.png"
$ws.Cells.Item(20, 13).Value = 86.4
$ws.Cells.Item(20, 14).Value = 67.8
$ws.Cells.Item(20, 15).Value = 77.1
$ws.Cells.Item(20, 16).Value = 35
$ws.Cells.Item(20, 17).Value = 9
$ws.Cells.Item(20, 18).Value = 9
$ws.Cells.Item(20, 19).Value = 9

# Row 21
$ws.Cells.Item(21, 3).Value = 2
$ws.Cells.Item(21, 6).Value = 87
$ws.Cells.Item(21, 12).Value = "stimuli/img_6nbgt.png"
$ws.Cells.Item(21, 13).Value = 78.45161290322581
$ws.Cells.Item(21, 14).Value = 57.83870967741935
$ws.Cells.Item(21, 15).Value = 68.14516129032258
$ws.Cells.Item(21, 17).Value = 7
$ws.Cells.Item(21, 18).Value = 7
$ws.Cells.Item(21, 19).Value = 7

# Row 22
$ws.Cells.Item(22, 3).Value = 2
$ws.Cells.Item(22, 6).Value = 88
$ws.Cells.Item(22, 12).Value = "stimuli/img_7pgd2.png"
$ws.Cells.Item(22, 13).Value = 78.59375
$ws.Cells.Item(22, 14).Value = 57.84375
$ws.Cells.Item(22, 15).Value = 68.21875
$ws.Cells.Item(22, 16).Value = 32
$ws.Cells.Item(22, 17).Value = 8
$ws.Cells.Item(22, 18).Value = 7
$ws.Cells.Item(22, 19).Value = 7

# Row 23
$ws.Cells.Item(23, 3).Value = 2
$ws.Cells.Item(23, 6).Value = 89
$ws.Cells.Item(23, 12).Value = "stimuli/img_84s7n.png"
$ws.Cells.Item(23, 13).Value = 11.03125
$ws.Cells.Item(23, 14).Value = 2.90625
$ws.Cells.Item(23, 15).Value = 6.96875
$ws.Cells.Item(23, 16).Value = 32
$ws.Cells.Item(23, 17).Value = 1
$ws.Cells.Item(23, 18).Value = 1
$ws.Cells.Item(23, 19).Value = 1

# Row 24
$ws.Cells.Item(24, 3).Value = 2
$ws.Cells.Item(24, 6).Value = 90
$ws.Cells.Item(24, 12).Value = "stimuli/img_24rt2.png"
$ws.Cells.Item(24, 13).Value = 55.26829268292683
$ws.Cells.Item(24, 14).Value = 34.19512195121951
$ws.Cells.Item(24, 15).Value = 44.73170731707317
$ws.Cells.Item(24, 16).Value = 41
$ws.Cells.Item(24, 17).Value = 3
$ws.Cells.Item(24, 18).Value = 3
$ws.Cells.Item(24, 19).Value = 3

# Row 25
$ws.Cells.Item(25, 3).Value = 2
$ws.Cells.Item(25, 6).Value = 91
$ws.Cells.Item(25, 12).Value = "stimuli/img_s9are.png"
$ws.Cells.Item(25, 13).Value = 90.14285714285714
$ws.Cells.Item(25, 14).Value = 75.22857142857143
$ws.Cells.Item(25, 15).Value = 82.68571428571428
$ws.Cells.Item(25, 16).Value = 35
$ws.Cells.Item(25, 17).Value = 10
$ws.Cells.Item(25, 18).Value = 10
$ws.Cells.Item(25, 19).Value = 10

# Row 26
$ws.Cells.Item(26, 3).Value = 2
$ws.Cells.Item(26, 6).Value = 92
$ws.Cells.Item(26, 8).Value = "kitchens"
$ws.Cells.Item(26, 9).Value = "target"
$ws.Cells.Item(26, 11).Value = "j"
$ws.Cells.Item(26, 12).Value = "stimuli/img_05flq.png"
$ws.Cells.Item(26, 13).Value = 47.10344827586207
$ws.Cells.Item(26, 14).Value = 25.72413793103448
$ws.Cells.Item(26, 15).Value = 36.41379310344828
$ws.Cells.Item(26, 16).Value = 29
$ws.Cells.Item(26, 17).Value = 1
$ws.Cells.Item(26, 18).Value = 1
$ws.Cells.Item(26, 19).Value = 1

# Row 27
$ws.Cells.Item(27, 3).Value = 2
$ws.Cells.Item(27, 6).Value = 93
$ws.Cells.Item(27, 8).Value = "bedrooms"
$ws.Cells.Item(27, 9).Value = "distractor"
$ws.Cells.Item(27, 11).Value = "f"
$ws.Cells.Item(27, 12).Value = "stimuli/img_scrdm.png"
$ws.Cells.Item(27, 13).Value = 78.675
$ws.Cells.Item(27, 14).Value = 57.9
$ws.Cells.Item(27, 15).Value = 68.2875
$ws.Cells.Item(27, 16).Value = 40
$ws.Cells.Item(27, 17).Value = 7
$ws.Cells.Item(27, 18).Value = 7
$ws.Cells.Item(27, 19).Value = 7
